$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 37.05583833333333
$ws.Range("H2").Value = 111.167515
$ws.Range("I2").Value = 0.008431126118266585
$ws.Range("J2").Value = 0.008431126118266585
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1788416666666666
$ws.Range("N2").Value = 0.5365249999999999
$ws.Range("O2").Value = 0.005632743940253072
$ws.Range("P2").Value = 0.005632743940253071
$ws.Range("Q2").Value = 6.627127887263888
$ws.Range("R2").Value = 59.64415098537498
$ws.Range("S2").Value = 0.0000474903745521755
$ws.Range("T2").Value = 0.0000474903745521755

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.05583833333333
$ws.Range("H3").Value = 111.167515
$ws.Range("I3").Value = 0.008431126118266585
$ws.Range("J3").Value = 0.008431126118266585
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.1427333333333334
$ws.Range("N3").Value = 0.4282
$ws.Range("O3").Value = 0.004495486613329045
$ws.Range("P3").Value = 0.004495486613329045
$ws.Range("Q3").Value = 5.289103324777779
$ws.Range("R3").Value = 47.601929923
$ws.Range("S3").Value = 0.00003790201459995631
$ws.Range("T3").Value = 0.00003790201459995631

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 37.05583833333333
$ws.Range("H4").Value = 111.167515
$ws.Range("I4").Value = 0.008431126118266585
$ws.Range("J4").Value = 0.008431126118266585
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.674834333333334
$ws.Range("N4").Value = 23.024503
$ws.Range("O4").Value = 0.2417242994279646
$ws.Range("P4").Value = 0.2417242994279646
$ws.Range("Q4").Value = 284.3974202911161
$ws.Range("R4").Value = 2559.576782620045
$ws.Range("S4").Value = 0.002038008054326805
$ws.Range("T4").Value = 0.002038008054326804

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.05583833333333
$ws.Range("H5").Value = 111.167515
$ws.Range("I5").Value = 0.008431126118266585
$ws.Range("J5").Value = 0.008431126118266585
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 23.42985866666666
$ws.Range("N5").Value = 70.289576
$ws.Range("O5").Value = 0.7379398597958301
$ws.Range("P5").Value = 0.7379398597958301
$ws.Range("Q5").Value = 868.2130549248488
$ws.Range("R5").Value = 7813.917494323639
$ws.Range("S5").Value = 0.006221664025634605
$ws.Range("T5").Value = 0.006221664025634605

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.05583833333333
$ws.Range("H6").Value = 111.167515
$ws.Range("I6").Value = 0.008431126118266585
$ws.Range("J6").Value = 0.008431126118266585
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3240953333333333
$ws.Range("N6").Value = 0.972286
$ws.Range("O6").Value = 0.01020761022262317
$ws.Range("P6").Value = 0.01020761022262317
$ws.Range("Q6").Value = 12.00962427658778
$ws.Range("R6").Value = 108.08661848929
$ws.Range("S6").Value = 0.00008606164915304322
$ws.Range("T6").Value = 0.00008606164915304322

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 50.89916233333333
$ws.Range("H7").Value = 152.697487
$ws.Range("I7").Value = 0.01158082710438721
$ws.Range("J7").Value = 0.01158082710438721
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.1788416666666666
$ws.Range("N7").Value = 0.5365249999999999
$ws.Range("O7").Value = 0.005632743940253072
$ws.Range("P7").Value = 0.005632743940253071
$ws.Range("Q7").Value = 9.102891023630553
$ws.Range("R7").Value = 81.92601921267499
$ws.Range("S7").Value = 0.0000652318336953556
$ws.Range("T7").Value = 0.0000652318336953556

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.89916233333333
$ws.Range("H8").Value = 152.697487
$ws.Range("I8").Value = 0.01158082710438721
$ws.Range("J8").Value = 0.01158082710438721
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1427333333333334
$ws.Range("N8").Value = 0.4282
$ws.Range("O8").Value = 0.004495486613329045
$ws.Range("P8").Value = 0.004495486613329045
$ws.Range("Q8").Value = 7.265007103711111
$ws.Range("R8").Value = 65.3850639334
$ws.Range("S8").Value = 0.00005206145321905089
$ws.Range("T8").Value = 0.00005206145321905089

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.89916233333333
$ws.Range("H9").Value = 152.697487
$ws.Range("I9").Value = 0.01158082710438721
$ws.Range("J9").Value = 0.01158082710438721
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.674834333333334
$ws.Range("N9").Value = 23.024503
$ws.Range("O9").Value = 0.2417242994279646
$ws.Range("P9").Value = 0.2417242994279646
$ws.Range("Q9").Value = 390.6426386137734
$ws.Range("R9").Value = 3515.783747523961
$ws.Range("S9").Value = 0.002799367318604383
$ws.Range("T9").Value = 0.002799367318604383

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.89916233333333
$ws.Range("H10").Value = 152.697487
$ws.Range("I10").Value = 0.01158082710438721
$ws.Range("J10").Value = 0.01158082710438721
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.42985866666666
$ws.Range("N10").Value = 70.289576
$ws.Range("O10").Value = 0.7379398597958301
$ws.Range("P10").Value = 0.7379398597958301
$ws.Range("Q10").Value = 1192.560179721723
$ws.Range("R10").Value = 10733.04161749551
$ws.Range("S10").Value = 0.00854595392973125
$ws.Range("T10").Value = 0.00854595392973125

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 50.89916233333333
$ws.Range("H11").Value = 152.697487
$ws.Range("I11").Value = 0.01158082710438721
$ws.Range("J11").Value = 0.01158082710438721
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.3240953333333333
$ws.Range("N11").Value = 0.972286
$ws.Range("O11").Value = 0.01020761022262317
$ws.Range("P11").Value = 0.01020761022262317
$ws.Range("Q11").Value = 16.49618098280911
$ws.Range("R11").Value = 148.465628845282
$ws.Range("S11").Value = 0.0001182125691371745
$ws.Range("T11").Value = 0.0001182125691371745

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2007.446289
$ws.Range("H12").Value = 6022.338867
$ws.Range("I12").Value = 0.4567440273772037
$ws.Range("J12").Value = 0.4567440273772037
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1788416666666666
$ws.Range("N12").Value = 0.5365249999999999
$ws.Range("O12").Value = 0.005632743940253072
$ws.Range("P12").Value = 0.005632743940253071
$ws.Range("Q12").Value = 359.015040068575
$ws.Range("R12").Value = 3231.135360617175
$ws.Range("S12").Value = 0.002572722152455727
$ws.Range("T12").Value = 0.002572722152455726

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2007.446289
$ws.Range("H13").Value = 6022.338867
$ws.Range("I13").Value = 0.4567440273772037
$ws.Range("J13").Value = 0.4567440273772037
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.1427333333333334
$ws.Range("N13").Value = 0.4282
$ws.Range("O13").Value = 0.004495486613329045
$ws.Range("P13").Value = 0.004495486613329045
$ws.Range("Q13").Value = 286.5295003166001
$ws.Range("R13").Value = 2578.7655028494
$ws.Range("S13").Value = 0.002053286660792214
$ws.Range("T13").Value = 0.002053286660792214

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2007.446289
$ws.Range("H14").Value = 6022.338867
$ws.Range("I14").Value = 0.4567440273772037
$ws.Range("J14").Value = 0.4567440273772037
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 7.674834333333334
$ws.Range("N14").Value = 23.024503
$ws.Range("O14").Value = 0.2417242994279646
$ws.Range("P14").Value = 0.2417242994279646
$ws.Range("Q14").Value = 15406.81770113979
$ws.Range("R14").Value = 138661.3593102581
$ws.Range("S14").Value = 0.1104061300356616
$ws.Range("T14").Value = 0.1104061300356616

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2007.446289
$ws.Range("H15").Value = 6022.338867
$ws.Range("I15").Value = 0.4567440273772037
$ws.Range("J15").Value = 0.4567440273772037
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 23.42985866666666
$ws.Range("N15").Value = 70.289576
$ws.Range("O15").Value = 0.7379398597958301
$ws.Range("P15").Value = 0.7379398597958301
$ws.Range("Q15").Value = 47034.18283219448
$ws.Range("R15").Value = 423307.6454897504
$ws.Range("S15").Value = 0.3370496235253164
$ws.Range("T15").Value = 0.3370496235253164

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2007.446289
$ws.Range("H16").Value = 6022.338867
$ws.Range("I16").Value = 0.4567440273772037
$ws.Range("J16").Value = 0.4567440273772037
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.3240953333333333
$ws.Range("N16").Value = 0.972286
$ws.Range("O16").Value = 0.01020761022262317
$ws.Range("P16").Value = 0.01020761022262317
$ws.Range("Q16").Value = 650.603974182218
$ws.Range("R16").Value = 5855.435767639962
$ws.Range("S16").Value = 0.004662265002977623
$ws.Range("T16").Value = 0.004662265002977623

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2293.273345666667
$ws.Range("H17").Value = 6879.820037
$ws.Range("I17").Value = 0.5217768014597114
$ws.Range("J17").Value = 0.5217768014597114
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.1788416666666666
$ws.Range("N17").Value = 0.5365249999999999
$ws.Range("O17").Value = 0.005632743940253072
$ws.Range("P17").Value = 0.005632743940253071
$ws.Range("Q17").Value = 410.1328272612695
$ws.Range("R17").Value = 3691.195445351424
$ws.Range("S17").Value = 0.002939035116586819
$ws.Range("T17").Value = 0.002939035116586819

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 2293.273345666667
$ws.Range("H18").Value = 6879.820037
$ws.Range("I18").Value = 0.5217768014597114
$ws.Range("J18").Value = 0.5217768014597114
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.1427333333333334
$ws.Range("N18").Value = 0.4282
$ws.Range("O18").Value = 0.004495486613329045
$ws.Range("P18").Value = 0.004495486613329045
$ws.Range("Q18").Value = 327.326548871489
$ws.Range("R18").Value = 2945.9389398434
$ws.Range("S18").Value = 0.002345640626107779
$ws.Range("T18").Value = 0.002345640626107779

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 2293.273345666667
$ws.Range("H19").Value = 6879.820037
$ws.Range("I19").Value = 0.5217768014597114
$ws.Range("J19").Value = 0.5217768014597114
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 7.674834333333334
$ws.Range("N19").Value = 23.024503
$ws.Range("O19").Value = 0.2417242994279646
$ws.Range("P19").Value = 0.2417242994279646
$ws.Range("Q19").Value = 17600.49300904074
$ws.Range("R19").Value = 158404.4370813666
$ws.Range("S19").Value = 0.1261261317906129
$ws.Range("T19").Value = 0.1261261317906129

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 2293.273345666667
$ws.Range("H20").Value = 6879.820037
$ws.Range("I20").Value = 0.5217768014597114
$ws.Range("J20").Value = 0.5217768014597114
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 23.42985866666666
$ws.Range("N20").Value = 70.289576
$ws.Range("O20").Value = 0.7379398597958301
$ws.Range("P20").Value = 0.7379398597958301
$ws.Range("Q20").Value = 53731.07037300381
$ws.Range("R20").Value = 483579.6333570343
$ws.Range("S20").Value = 0.3850398997138961
$ws.Range("T20").Value = 0.3850398997138961

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 2293.273345666667
$ws.Range("H21").Value = 6879.820037
$ws.Range("I21").Value = 0.5217768014597114
$ws.Range("J21").Value = 0.5217768014597114
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.3240953333333333
$ws.Range("N21").Value = 0.972286
$ws.Range("O21").Value = 0.01020761022262317
$ws.Range("P21").Value = 0.01020761022262317
$ws.Range("Q21").Value = 743.2391893882871
$ws.Range("R21").Value = 6689.152704494582
$ws.Range("S21").Value = 0.005326094212507773
$ws.Range("T21").Value = 0.005326094212507773

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 6.448603666666666
$ws.Range("H22").Value = 19.345811
$ws.Range("I22").Value = 0.00146721794043115
$ws.Range("J22").Value = 0.00146721794043115
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 0.1788416666666666
$ws.Range("N22").Value = 0.5365249999999999
$ws.Range("O22").Value = 0.005632743940253072
$ws.Range("P22").Value = 0.005632743940253071
$ws.Range("Q22").Value = 1.153279027419444
$ws.Range("R22").Value = 10.379511246775
$ws.Range("S22").Value = 0.000008264462962994153
$ws.Range("T22").Value = 0.000008264462962994151

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 6.448603666666666
$ws.Range("H23").Value = 19.345811
$ws.Range("I23").Value = 0.00146721794043115
$ws.Range("J23").Value = 0.00146721794043115
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.1427333333333334
$ws.Range("N23").Value = 0.4282
$ws.Range("O23").Value = 0.004495486613329045
$ws.Range("P23").Value = 0.004495486613329045
$ws.Range("Q23").Value = 0.9204306966888889
$ws.Range("R23").Value = 8.283876270199999
$ws.Range("S23").Value = 0.000006595858610044447
$ws.Range("T23").Value = 0.000006595858610044447

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 6.448603666666666
$ws.Range("H24").Value = 19.345811
$ws.Range("I24").Value = 0.00146721794043115
$ws.Range("J24").Value = 0.00146721794043115
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 7.674834333333334
$ws.Range("N24").Value = 23.024503
$ws.Range("O24").Value = 0.2417242994279646
$ws.Range("P24").Value = 0.2417242994279646
$ws.Range("Q24").Value = 49.49196482299256
$ws.Range("R24").Value = 445.427683406933
$ws.Range("S24").Value = 0.0003546622287588608
$ws.Range("T24").Value = 0.0003546622287588608

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 6.448603666666666
$ws.Range("H25").Value = 19.345811
$ws.Range("I25").Value = 0.00146721794043115
$ws.Range("J25").Value = 0.00146721794043115
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 23.42985866666666
$ws.Range("N25").Value = 70.289576
$ws.Range("O25").Value = 0.7379398597958301
$ws.Range("P25").Value = 0.7379398597958301
$ws.Range("Q25").Value = 151.0898725073484
$ws.Range("R25").Value = 1359.808852566136
$ws.Range("S25").Value = 0.00108271860125169
$ws.Range("T25").Value = 0.00108271860125169

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 6.448603666666666
$ws.Range("H26").Value = 19.345811
$ws.Range("I26").Value = 0.00146721794043115
$ws.Range("J26").Value = 0.00146721794043115
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 0.3240953333333333
$ws.Range("N26").Value = 0.972286
$ws.Range("O26").Value = 0.01020761022262317
$ws.Range("P26").Value = 0.01020761022262317
$ws.Range("Q26").Value = 2.089962354882889
$ws.Range("R26").Value = 18.809661193946
$ws.Range("S26").Value = 0.00001497678884756113
$ws.Range("T26").Value = 0.00001497678884756113

